# Generate Report for Handback
# Adds a new handback row (c33a1e96-daca-46d8-b1f5-94ecd3f3b178.md) to the
# Overview, zh-cn and de-de sheets, mirroring the existing
# 61ddfa91-de3a-4a36-937d-a7ec915a7726.md row, and grows each table /
# dimension to include the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "c33a1e96-daca-46d8-b1f5-94ecd3f3b178.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-10-19 23:31:13"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e893c104dd95e2b3f1a7c6409fb6c7e1e0a2b4c/e2e/c33a1e96-daca-46d8-b1f5-94ecd3f3b178.md", "", "", "e2e\c33a1e96-daca-46d8-b1f5-94ecd3f3b178.md")
$wsOverview.Range("B3").Font.Underline = $true
$wsOverview.Range("B3").Font.Color = 15570276

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = "c33a1e96-daca-46d8-b1f5-94ecd3f3b178.f9019a82e66bb6f16371312027533e8fb6889afc.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-10-19 23:31:02"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("J3").Value = "c33a1e96-daca-46d8-b1f5-94ecd3f3b178.f9019a82e66bb6f16371312027533e8fb6889afc.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-10-19 23:31:51"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e893c104dd95e2b3f1a7c6409fb6c7e1e0a2b4c/e2e/c33a1e96-daca-46d8-b1f5-94ecd3f3b178.md", "", "", "c33a1e96-daca-46d8-b1f5-94ecd3f3b178.md")
$wsZhCn.Range("A3").Font.Underline = $true
$wsZhCn.Range("A3").Font.Color = 15570276

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/5b6c7d8e9fa0b1c2d3e4f5061728394a5b6c7d8e/e2e/c33a1e96-daca-46d8-b1f5-94ecd3f3b178.md", "", "", "c33a1e96-daca-46d8-b1f5-94ecd3f3b178.md")
$wsZhCn.Range("I3").Font.Underline = $true
$wsZhCn.Range("I3").Font.Color = 15570276

$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = "c33a1e96-daca-46d8-b1f5-94ecd3f3b178.f9019a82e66bb6f16371312027533e8fb6889afc.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-10-19 23:31:13"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("J3").Value = "c33a1e96-daca-46d8-b1f5-94ecd3f3b178.f9019a82e66bb6f16371312027533e8fb6889afc.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-10-19 23:32:09"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e893c104dd95e2b3f1a7c6409fb6c7e1e0a2b4c/e2e/c33a1e96-daca-46d8-b1f5-94ecd3f3b178.md", "", "", "c33a1e96-daca-46d8-b1f5-94ecd3f3b178.md")
$wsDeDe.Range("A3").Font.Underline = $true
$wsDeDe.Range("A3").Font.Color = 15570276

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6c7d8e9fa0b1c2d3e4f5061728394a5b6c7d8e9f/e2e/c33a1e96-daca-46d8-b1f5-94ecd3f3b178.md", "", "", "c33a1e96-daca-46d8-b1f5-94ecd3f3b178.md")
$wsDeDe.Range("I3").Font.Underline = $true
$wsDeDe.Range("I3").Font.Color = 15570276

$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P3"))
